# Reorder the "Recorded By" (column G) author lists on the
# "Session Analysis Results" sheet so "System" no longer sorts first.
#
# Mapping applied (exact-value swap, everything else left untouched):
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "admin@admin.com, dnasr281@gmail.com" -> "dnasr281@gmail.com, admin@admin.com"
#   "backup@backdoor.com, System, system" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
    "backup@backdoor.com, System, system" = "system, backup@backdoor.com, System"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2
    if ($null -ne $value -and $map.ContainsKey($value)) {
        $cell.Value = $map[$value]
    }
}
